$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had columns: A Sites | B Nombre d'individus observées | C Nombre d'espèces
# observées | D N0 (rarefied) | E N1 (rarefied) | F N2 (rarefied) | G Type
#
# This commit swaps in non-rarefied Hill numbers: the old rarefied "N0" column
# (D) is dropped entirely (species richness, column C, now doubles as N0), and
# N1/N2 are recomputed from the raw, non-rarefied abundances. Deleting column D
# shifts the old N1/N2/Type columns left into D/E/F, which already gives D/E/F
# the right headers ("N1"/"N2"/"Type") for free.
$ws.Columns.Item(4).Delete()

# Fix up the two header cells that still need new text.
$ws.Range("B1").Value = "abon"
$ws.Range("C1").Value = "N0"

# Recalculated (non-rarefied) N1 / N2 values per site.
$n1 = @{2=19; 3=13; 4=14; 5=22; 6=12; 7=17; 8=29; 9=21; 10=15; 11=19; 12=22; 13=20; 14=21; 15=24; 16=18}
$n2 = @{2=13; 3=7; 4=7; 5=17; 6=7; 7=9; 8=21; 9=10; 10=8; 11=11; 12=14; 13=10; 14=13; 15=15; 16=10}

for ($row = 2; $row -le 16; $row++) {
    $ws.Cells.Item($row, 4).Value = $n1[$row]
    $ws.Cells.Item($row, 5).Value = $n2[$row]
}

# "Rue de l'Egalite" loses one observed individual, which also happened to be
# the sole record of one species, so both its individual and species counts
# drop by one; the grand total of individuals drops by one accordingly.
$ws.Range("B10").Value = 121
$ws.Range("C10").Value = 31
$ws.Range("B17").Value = 1902
